# Generate Report for Handoff
# Renames the backing e2e markdown test file from the old GUID-based name to a new one,
# refreshes the generated xliff file names / timestamps, and clears out the
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime" columns
# on the per-locale sheets now that a fresh handoff cycle has started.

$wb = $excel.ActiveWorkbook

# New e2e backing file (old one was 849fcfc5-ddc8-4b60-9160-1e3a8f26d5bb.md)
$newGuidFile = "fa5cde86-5f3a-49e2-985f-dafd8cb5bf0c.md"
$newGuidPath = "e2e\fa5cde86-5f3a-49e2-985f-dafd8cb5bf0c.md"

# New generated xliff file names for this handoff cycle
$newZhXlf = "fa5cde86-5f3a-49e2-985f-dafd8cb5bf0c.a614736447dd9fa95efbd58a20b0c48f8ff778b0.zh-cn.xlf"
$newDeXlf = "fa5cde86-5f3a-49e2-985f-dafd8cb5bf0c.a614736447dd9fa95efbd58a20b0c48f8ff778b0.de-de.xlf"

function Remove-HyperlinkAt($ws, $addr) {
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
}

function Set-HyperlinkDisplay($ws, $addr, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet: only the file name / path text (and the matching hyperlink
# display text) change; the generate date is refreshed as well.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newGuidFile
$wsOverview.Range("B2").Value = $newGuidPath
Set-HyperlinkDisplay $wsOverview '$B$2' $newGuidPath
$wsOverview.Range("G2").Value = "2016-09-05 07:09:49"
# Target stored width is 39.7183053152902; Excel can only persist column widths at
# whole-pixel granularity, so pick the ColumnWidth input that lands on the closest
# achievable stored width (39.6666666666667).
$wsOverview.Columns.Item(1).ColumnWidth = 38.8333333333333

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newGuidFile
Set-HyperlinkDisplay $wsZh '$A$2' $newGuidFile

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-09-05 07:09:44"

Remove-HyperlinkAt $wsZh '$I$2'
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"

$wsZh.Range("J2").Value = ""

$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# Target stored widths are 39.7183053152902 / 18.6506053379604 / 21.7054770333426;
# closest achievable (whole-pixel) stored widths are 39.6666666666667 / 18.6666666666667 / 21.6666666666667.
$wsZh.Columns.Item(1).ColumnWidth = 38.8333333333333
$wsZh.Columns.Item(9).ColumnWidth = 17.8333333333333
$wsZh.Columns.Item(10).ColumnWidth = 20.8333333333333

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newGuidFile
Set-HyperlinkDisplay $wsDe '$A$2' $newGuidFile

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = "2016-09-05 07:09:49"

Remove-HyperlinkAt $wsDe '$I$2'
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"

$wsDe.Range("J2").Value = ""

$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

# Target stored widths are 39.7183053152902 / 18.6506053379604 / 21.7054770333426;
# closest achievable (whole-pixel) stored widths are 39.6666666666667 / 18.6666666666667 / 21.6666666666667.
$wsDe.Columns.Item(1).ColumnWidth = 38.8333333333333
$wsDe.Columns.Item(9).ColumnWidth = 17.8333333333333
$wsDe.Columns.Item(10).ColumnWidth = 20.8333333333333

$wb.Save()
